$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $cell = $ws.Range($cellAddr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue "D2" "53.770.65"
Set-TextValue "E2" "  -4.77%  "

Set-TextValue "D3" "2.214.98"
Set-TextValue "E3" "  -6.73%  "

Set-TextValue "D4" "1.01"
Set-TextValue "E4" "  +0.48%  "

Set-TextValue "D5" "488.33"
Set-TextValue "E5" "  -3.78%  "

Set-TextValue "D6" "125.84"
Set-TextValue "E6" "  -3.61%  "

Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  +0.11%  "

Set-TextValue "D8" "0.524"
Set-TextValue "E8" "  -3.86%  "

Set-TextValue "D9" "2.244.54"
Set-TextValue "E9" "  -5.85%  "

Set-TextValue "D10" "0.0924"
Set-TextValue "E10" "  -6.32%  "

Set-TextValue "E11" "  -0.21%  "

Set-TextValue "D12" "0.321"
Set-TextValue "E12" "  -3.15%  "

Set-TextValue "E13" "  -4.55%  "

Set-TextValue "D14" "2.614.93"
Set-TextValue "E14" "  -6.44%  "

Set-TextValue "D15" "21.29"
Set-TextValue "E15" "  -1.44%  "

Set-TextValue "D16" "53.770.44"
Set-TextValue "E16" "  -4.70%  "

Set-TextValue "E17" "  -3.68%  "

Set-TextValue "D18" "2.229.22"
Set-TextValue "E18" "  -6.89%  "

Set-TextValue "B19" "Polkadot"
Set-TextValue "C19" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D19" "4.00"
Set-TextValue "E19" "  -1.24%  "

Set-TextValue "B20" "Chainlink"
Set-TextValue "C20" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "9.66"
Set-TextValue "E20" "  -3.92%  "

Set-TextValue "D21" "296.90"
Set-TextValue "E21" "  -4.06%  "

Set-TextValue "D22" "6.26"
Set-TextValue "E22" "  -0.37%  "

Set-TextValue "D23" "0.999"
Set-TextValue "E23" "  -0.04%  "

Set-TextValue "D24" "63.95"
Set-TextValue "E24" "  -3.42%  "

Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  +0.32%  "

Set-TextValue "D26" "0.369"
Set-TextValue "E26" "  -0.77%  "

Set-TextValue "E27" "  -1.09%  "

Set-TextValue "D28" "2.314.48"

Set-TextValue "E29" "  -2.79%  "

Set-TextValue "D30" "163.09"
Set-TextValue "E30" "  -5.90%  "

Set-TextValue "E31" "  -3.85%  "

Set-TextValue "D33" "5.81"
Set-TextValue "E33" "  -0.68%  "

Set-TextValue "D34" "0.0₃0671"
Set-TextValue "E34" "  -5.80%  "

Set-TextValue "E36" "  -1.65%  "

Set-TextValue "D37" "17.32"
Set-TextValue "E37" "  -1.83%  "

Set-TextValue "E38" "  -0.38%  "

Set-TextValue "D39" "0.837"
Set-TextValue "E39" "  +1.46%  "

Set-TextValue "E40" "  -3.23%  "

Set-TextValue "E41" "  -3.44%  "

Set-TextValue "D42" "1.39"
Set-TextValue "E42" "  -0.26%  "

Set-TextValue "D43" "0.370"
Set-TextValue "E43" "  -0.52%  "

Set-TextValue "D44" "127.93"
Set-TextValue "E44" "  +1.49%  "

Set-TextValue "D45" "3.30"
Set-TextValue "E45" "  -2.26%  "

Set-TextValue "D46" "4.79"
Set-TextValue "E46" "  -3.41%  "

Set-TextValue "D47" "0.0887"
Set-TextValue "E47" "  -1.27%  "

Set-TextValue "D48" "0.538"
Set-TextValue "E48" "  -4.75%  "

Set-TextValue "D49" "237.99"
Set-TextValue "E49" "  -0.58%  "

Set-TextValue "E50" "  -1.89%  "

Set-TextValue "E51" "  -3.12%  "
